# Sync the cryptocurrency price/volume table to the latest scrape.
# Every write goes through Set-TextValue so values stay plain text
# (matching the source sheet's inline-string cells) instead of being
# auto-coerced to numbers by Excel's type inference (e.g. "17.86",
# "10.00", "22.00" would otherwise lose trailing zeros / become floats).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"   # drop the temporary text format again, keep default styling
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue "D2" "68.472.12"
Set-TextValue "E2" "  +0.53%  "

Set-TextValue "D3" "2.690.57"
Set-TextValue "E3" "  +1.83%  "

Set-TextValue "E4" "  +0.01%  "

Set-TextValue "D5" "598.46"
Set-TextValue "E5" "  +0.27%  "

Set-TextValue "D6" "159.83"
Set-TextValue "E6" "  +2.79%  "

Set-TextValue "E7" "  +0.02%  "

Set-TextValue "E8" "  +0.39%  "

Set-TextValue "D9" "2.689.97"
Set-TextValue "E9" "  +1.85%  "

Set-TextValue "E10" "  +0.26%  "

Set-TextValue "E11" "  -0.61%  "

Set-TextValue "E12" "  +1.17%  "

Set-TextValue "E13" "  +2.55%  "

Set-TextValue "D14" "28.21"
Set-TextValue "E14" "  +1.06%  "

Set-TextValue "D15" "3.179.31"
Set-TextValue "E15" "  +1.78%  "

Set-TextValue "E16" "  -0.14%  "

Set-TextValue "D17" "68.393.13"
Set-TextValue "E17" "  +0.47%  "

Set-TextValue "D18" "2.690.72"
Set-TextValue "E18" "  +0.99%  "

Set-TextValue "D19" "11.85"
Set-TextValue "E19" "  +4.61%  "

Set-TextValue "D20" "366.86"
Set-TextValue "E20" "  +1.09%  "

Set-TextValue "E21" "  +3.25%  "

Set-TextValue "D22" "4.54"
Set-TextValue "E22" "  +3.09%  "

Set-TextValue "E23" "  +2.01%  "

Set-TextValue "D24" "2.12"
Set-TextValue "E24" "  +2.76%  "

Set-TextValue "D25" "74.50"
Set-TextValue "E25" "  -0.12%  "

Set-TextValue "D26" "0.998"
Set-TextValue "E26" "  -0.15%  "

Set-TextValue "D27" "10.00"
Set-TextValue "E27" "  +2.88%  "

Set-TextValue "E29" "  +0.46%  "

Set-TextValue "D30" "1.00"
Set-TextValue "E30" "  +0.21%  "

Set-TextValue "D31" "571.52"
Set-TextValue "E31" "  +3.11%  "

Set-TextValue "E33" "  +3.63%  "

Set-TextValue "E34" "  +5.52%  "

Set-TextValue "E35" "  +2.67%  "

Set-TextValue "E36" "  +6.65%  "

Set-TextValue "E37" "  -0.01%  "

Set-TextValue "D38" "161.82"
Set-TextValue "E38" "  +0.41%  "

Set-TextValue "D39" "19.85"
Set-TextValue "E39" "  +2.24%  "

Set-TextValue "D40" "0.379"
Set-TextValue "E40" "  +1.97%  "

Set-TextValue "E41" "  +2.02%  "

Set-TextValue "E42" "  +1.59%  "

# Row 43/44: dogwifhat and WhiteBITCoin swap rank positions
Set-TextValue "B43" "dogwifhat"
Set-TextValue "C43" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D43" "2.66"
Set-TextValue "E43" "  +1.94%  "

Set-TextValue "B44" "WhiteBITCoin"
Set-TextValue "C44" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D44" "17.86"
Set-TextValue "E44" "  +0.39%  "

Set-TextValue "E45" "  +0.05%  "

Set-TextValue "E46" "  -6.22%  "

Set-TextValue "D47" "157.23"
Set-TextValue "E47" "  -1.35%  "

Set-TextValue "E48" "  +7.32%  "

Set-TextValue "D49" "1.77"
Set-TextValue "E49" "  +4.86%  "

Set-TextValue "D50" "0.597"
Set-TextValue "E50" "  +6.34%  "

Set-TextValue "D51" "22.00"
Set-TextValue "E51" "  +0.06%  "
